$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-07-09 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-10 Monday", 2) | Out-Null
$d.Content.Find.Execute("99-57=42", $true, $false, $false, $false, $false, $true, 1, $false, "78+2=80", 2) | Out-Null
$d.Content.Find.Execute("17+81=98", $true, $false, $false, $false, $false, $true, 1, $false, "98-63=35", 2) | Out-Null
$d.Content.Find.Execute("53-49=4", $true, $false, $false, $false, $false, $true, 1, $false, "3+57=60", 2) | Out-Null
$d.Content.Find.Execute("61-10=51", $true, $false, $false, $false, $false, $true, 1, $false, "44+15=59", 2) | Out-Null
$d.Content.Find.Execute("54+33=87", $true, $false, $false, $false, $false, $true, 1, $false, "73-38=35", 2) | Out-Null
$d.Content.Find.Execute("45-19=26", $true, $false, $false, $false, $false, $true, 1, $false, "46+8=54", 2) | Out-Null
$d.Content.Find.Execute("62+8=70", $true, $false, $false, $false, $false, $true, 1, $false, "62+18=80", 2) | Out-Null
$d.Content.Find.Execute("55-19=36", $true, $false, $false, $false, $false, $true, 1, $false, "63-60=3", 2) | Out-Null
$d.Content.Find.Execute("71+12=83", $true, $false, $false, $false, $false, $true, 1, $false, "29+7=36", 2) | Out-Null
$d.Content.Find.Execute("33+16=49", $true, $false, $false, $false, $false, $true, 1, $false, "2+27=29", 2) | Out-Null
$d.Content.Find.Execute("56-50=6", $true, $false, $false, $false, $false, $true, 1, $false, "33+64=97", 2) | Out-Null
$d.Content.Find.Execute("19+18=37", $true, $false, $false, $false, $false, $true, 1, $false, "54-15=39", 2) | Out-Null
$d.Content.Find.Execute("87-38=49", $true, $false, $false, $false, $false, $true, 1, $false, "5+84=89", 2) | Out-Null
$d.Content.Find.Execute("91+6=97", $true, $false, $false, $false, $false, $true, 1, $false, "26+35=61", 2) | Out-Null
$d.Content.Find.Execute("89-84=5", $true, $false, $false, $false, $false, $true, 1, $false, "41+20=61", 2) | Out-Null
$d.Content.Find.Execute("31+56=87", $true, $false, $false, $false, $false, $true, 1, $false, "14+6=20", 2) | Out-Null
$d.Content.Find.Execute("84-2=82", $true, $false, $false, $false, $false, $true, 1, $false, "92+0=92", 2) | Out-Null
$d.Content.Find.Execute("58-29=29", $true, $false, $false, $false, $false, $true, 1, $false, "35+13=48", 2) | Out-Null
$d.Content.Find.Execute("31+61=92", $true, $false, $false, $false, $false, $true, 1, $false, "52-1=51", 2) | Out-Null
$d.Content.Find.Execute("37+0=37", $true, $false, $false, $false, $false, $true, 1, $false, "99-62=37", 2) | Out-Null
$d.Content.Find.Execute("96-53=43", $true, $false, $false, $false, $false, $true, 1, $false, "23+13=36", 2) | Out-Null
$d.Content.Find.Execute("34+12=46", $true, $false, $false, $false, $false, $true, 1, $false, "24-8=16", 2) | Out-Null
$d.Content.Find.Execute("62-27=35", $true, $false, $false, $false, $false, $true, 1, $false, "77-43=34", 2) | Out-Null
$d.Content.Find.Execute("61+6=67", $true, $false, $false, $false, $false, $true, 1, $false, "70+19=89", 2) | Out-Null
$d.Content.Find.Execute("78-75=3", $true, $false, $false, $false, $false, $true, 1, $false, "94-5=89", 2) | Out-Null
$d.Content.Find.Execute("30+68=98", $true, $false, $false, $false, $false, $true, 1, $false, "87-59=28", 2) | Out-Null
$d.Content.Find.Execute("41-30=11", $true, $false, $false, $false, $false, $true, 1, $false, "50+22=72", 2) | Out-Null
$d.Content.Find.Execute("85-10=75", $true, $false, $false, $false, $false, $true, 1, $false, "63+12=75", 2) | Out-Null
$d.Content.Find.Execute("25+18=43", $true, $false, $false, $false, $false, $true, 1, $false, "13+76=89", 2) | Out-Null
$d.Content.Find.Execute("30+19=49", $true, $false, $false, $false, $false, $true, 1, $false, "28-19=9", 2) | Out-Null
$d.Content.Find.Execute("56+18=74", $true, $false, $false, $false, $false, $true, 1, $false, "64-9=55", 2) | Out-Null
$d.Content.Find.Execute("32-19=13", $true, $false, $false, $false, $false, $true, 1, $false, "2+33=35", 2) | Out-Null
$d.Content.Find.Execute("33-30=3", $true, $false, $false, $false, $false, $true, 1, $false, "21-1=20", 2) | Out-Null
$d.Content.Find.Execute("91-33=58", $true, $false, $false, $false, $false, $true, 1, $false, "31+59=90", 2) | Out-Null
$d.Content.Find.Execute("87-84=3", $true, $false, $false, $false, $false, $true, 1, $false, "69-34=35", 2) | Out-Null
$d.Content.Find.Execute("43+30=73", $true, $false, $false, $false, $false, $true, 1, $false, "33-31=2", 2) | Out-Null
$d.Content.Find.Execute("10+46=56", $true, $false, $false, $false, $false, $true, 1, $false, "22+39=61", 2) | Out-Null
$d.Content.Find.Execute("11+36=47", $true, $false, $false, $false, $false, $true, 1, $false, "19+10=29", 2) | Out-Null
$d.Content.Find.Execute("89-3=86", $true, $false, $false, $false, $false, $true, 1, $false, "71-64=7", 2) | Out-Null
$d.Content.Find.Execute("46+20=66", $true, $false, $false, $false, $false, $true, 1, $false, "76-23=53", 2) | Out-Null
$d.Content.Find.Execute("54+35=89", $true, $false, $false, $false, $false, $true, 1, $false, "90-2=88", 2) | Out-Null
$d.Content.Find.Execute("23-1=22", $true, $false, $false, $false, $false, $true, 1, $false, "14+11=25", 2) | Out-Null
$d.Content.Find.Execute("66+17=83", $true, $false, $false, $false, $false, $true, 1, $false, "94-18=76", 2) | Out-Null
$d.Content.Find.Execute("55-41=14", $true, $false, $false, $false, $false, $true, 1, $false, "49+39=88", 2) | Out-Null
$d.Content.Find.Execute("68-15=53", $true, $false, $false, $false, $false, $true, 1, $false, "61+37=98", 2) | Out-Null
$d.Content.Find.Execute("49+34=83", $true, $false, $false, $false, $false, $true, 1, $false, "10+29=39", 2) | Out-Null
$d.Content.Find.Execute("54-53=1", $true, $false, $false, $false, $false, $true, 1, $false, "71+6=77", 2) | Out-Null
$d.Content.Find.Execute("75-14=61", $true, $false, $false, $false, $false, $true, 1, $false, "96-62=34", 2) | Out-Null
$d.Content.Find.Execute("43+51=94", $true, $false, $false, $false, $false, $true, 1, $false, "71-66=5", 2) | Out-Null
$d.Content.Find.Execute("56-19=37", $true, $false, $false, $false, $false, $true, 1, $false, "89-67=22", 2) | Out-Null
$d.Content.Find.Execute("43-35=8", $true, $false, $false, $false, $false, $true, 1, $false, "72-16=56", 2) | Out-Null
$d.Content.Find.Execute("70-22=48", $true, $false, $false, $false, $false, $true, 1, $false, "40-4=36", 2) | Out-Null
$d.Content.Find.Execute("6+40=46", $true, $false, $false, $false, $false, $true, 1, $false, "8+61=69", 2) | Out-Null
$d.Content.Find.Execute("71-22=49", $true, $false, $false, $false, $false, $true, 1, $false, "20+71=91", 2) | Out-Null
$d.Content.Find.Execute("80-62=18", $true, $false, $false, $false, $false, $true, 1, $false, "45-16=29", 2) | Out-Null
$d.Content.Find.Execute("94-24=70", $true, $false, $false, $false, $false, $true, 1, $false, "30+58=88", 2) | Out-Null
$d.Content.Find.Execute("84-56=28", $true, $false, $false, $false, $false, $true, 1, $false, "60-54=6", 2) | Out-Null
$d.Content.Find.Execute("88-77=11", $true, $false, $false, $false, $false, $true, 1, $false, "11+25=36", 2) | Out-Null
$d.Content.Find.Execute("33+6=39", $true, $false, $false, $false, $false, $true, 1, $false, "47+13=60", 2) | Out-Null
$d.Content.Find.Execute("77-18=59", $true, $false, $false, $false, $false, $true, 1, $false, "33+1=34", 2) | Out-Null
$d.Content.Find.Execute("44-8=36", $true, $false, $false, $false, $false, $true, 1, $false, "12+83=95", 2) | Out-Null
$d.Content.Find.Execute("38+45=83", $true, $false, $false, $false, $false, $true, 1, $false, "14+19=33", 2) | Out-Null
$d.Content.Find.Execute("27-13=14", $true, $false, $false, $false, $false, $true, 1, $false, "95-0=95", 2) | Out-Null
$d.Content.Find.Execute("31+60=91", $true, $false, $false, $false, $false, $true, 1, $false, "75-40=35", 2) | Out-Null
$d.Content.Find.Execute("68-49=19", $true, $false, $false, $false, $false, $true, 1, $false, "57-28=29", 2) | Out-Null
$d.Content.Find.Execute("28+52=80", $true, $false, $false, $false, $false, $true, 1, $false, "77-23=54", 2) | Out-Null
$d.Content.Find.Execute("46+13=59", $true, $false, $false, $false, $false, $true, 1, $false, "50-12=38", 2) | Out-Null
$d.Content.Find.Execute("73+14=87", $true, $false, $false, $false, $false, $true, 1, $false, "8+55=63", 2) | Out-Null
$d.Content.Find.Execute("46-25=21", $true, $false, $false, $false, $false, $true, 1, $false, "98-14=84", 2) | Out-Null
$d.Content.Find.Execute("64-40=24", $true, $false, $false, $false, $false, $true, 1, $false, "43+8=51", 2) | Out-Null
$d.Content.Find.Execute("83-34=49", $true, $false, $false, $false, $false, $true, 1, $false, "1+1=2", 2) | Out-Null
$d.Content.Find.Execute("93-59=34", $true, $false, $false, $false, $false, $true, 1, $false, "21+78=99", 2) | Out-Null
$d.Content.Find.Execute("90-48=42", $true, $false, $false, $false, $false, $true, 1, $false, "46-3=43", 2) | Out-Null
$d.Content.Find.Execute("21+43=64", $true, $false, $false, $false, $false, $true, 1, $false, "0+34=34", 2) | Out-Null
$d.Content.Find.Execute("29+47=76", $true, $false, $false, $false, $false, $true, 1, $false, "59-23=36", 2) | Out-Null
$d.Content.Find.Execute("1+96=97", $true, $false, $false, $false, $false, $true, 1, $false, "26+6=32", 2) | Out-Null
$d.Content.Find.Execute("84-70=14", $true, $false, $false, $false, $false, $true, 1, $false, "91-77=14", 2) | Out-Null
$d.Content.Find.Execute("95-75=20", $true, $false, $false, $false, $false, $true, 1, $false, "17+59=76", 2) | Out-Null
$d.Content.Find.Execute("82-58=24", $true, $false, $false, $false, $false, $true, 1, $false, "50-15=35", 2) | Out-Null
$d.Content.Find.Execute("59-38=21", $true, $false, $false, $false, $false, $true, 1, $false, "8+10=18", 2) | Out-Null
$d.Content.Find.Execute("56+39=95", $true, $false, $false, $false, $false, $true, 1, $false, "39-16=23", 2) | Out-Null
$d.Content.Find.Execute("73-64=9", $true, $false, $false, $false, $false, $true, 1, $false, "64-13=51", 2) | Out-Null
$d.Content.Find.Execute("17+77=94", $true, $false, $false, $false, $false, $true, 1, $false, "36+4=40", 2) | Out-Null
$d.Content.Find.Execute("70-14=56", $true, $false, $false, $false, $false, $true, 1, $false, "3+65=68", 2) | Out-Null
$d.Content.Find.Execute("57-46=11", $true, $false, $false, $false, $false, $true, 1, $false, "91-72=19", 2) | Out-Null
$d.Content.Find.Execute("2+79=81", $true, $false, $false, $false, $false, $true, 1, $false, "56-37=19", 2) | Out-Null
$d.Content.Find.Execute("83+7=90", $true, $false, $false, $false, $false, $true, 1, $false, "86-75=11", 2) | Out-Null
$d.Content.Find.Execute("42+53=95", $true, $false, $false, $false, $false, $true, 1, $false, "96-12=84", 2) | Out-Null
$d.Content.Find.Execute("69-49=20", $true, $false, $false, $false, $false, $true, 1, $false, "33-28=5", 2) | Out-Null
$d.Content.Find.Execute("99-54=45", $true, $false, $false, $false, $false, $true, 1, $false, "31-24=7", 2) | Out-Null
$d.Content.Find.Execute("74-30=44", $true, $false, $false, $false, $false, $true, 1, $false, "76-72=4", 2) | Out-Null
$d.Content.Find.Execute("60+36=96", $true, $false, $false, $false, $false, $true, 1, $false, "89-2=87", 2) | Out-Null
$d.Content.Find.Execute("1+92=93", $true, $false, $false, $false, $false, $true, 1, $false, "26-3=23", 2) | Out-Null
$d.Content.Find.Execute("12+71=83", $true, $false, $false, $false, $false, $true, 1, $false, "28+7=35", 2) | Out-Null
$d.Content.Find.Execute("1+60=61", $true, $false, $false, $false, $false, $true, 1, $false, "0+52=52", 2) | Out-Null
$d.Content.Find.Execute("35+7=42", $true, $false, $false, $false, $false, $true, 1, $false, "30-19=11", 2) | Out-Null
$d.Content.Find.Execute("95-87=8", $true, $false, $false, $false, $false, $true, 1, $false, "79-22=57", 2) | Out-Null
$d.Content.Find.Execute("18+11=29", $true, $false, $false, $false, $false, $true, 1, $false, "87-65=22", 2) | Out-Null
$d.Content.Find.Execute("26+9=35", $true, $false, $false, $false, $false, $true, 1, $false, "41+53=94", 2) | Out-Null
$d.Content.Find.Execute("37+17=54", $true, $false, $false, $false, $false, $true, 1, $false, "10+70=80", 2) | Out-Null

Write-Output "Replacements applied: 101"
